# Examensarbete.xlsx edit: "Made some progress on converting the ica data."
#
# - Insert a new "Frontend" worksheet right after "Sheet1" (and before
#   "Ica-Api"), pre-populated with three notes about the frontend approach.
# - Move the active tab/selection onto the new Frontend sheet.
# - Nudge Sheet1's remembered selection.
# - Clear the "tabSelected" flag that used to sit on Willys-Api (it moves to
#   the new Frontend sheet automatically once that becomes the active sheet).

$wb = $excel.ActiveWorkbook

# --- Sheet1: update remembered selection, then make it active so later adds
#     anchor correctly relative to it -------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate() | Out-Null
$sheet1.Range("C5").Select() | Out-Null

# --- Willys-Api: drop its old "last active" selection; the new Frontend
#     sheet takes over as the active tab, so this sheet simply keeps its own
#     remembered cell selection but is no longer the active tab -----------
$willys = $wb.Worksheets.Item("Willys-Api")
$willys.Range("C6").Select() | Out-Null

# --- New "Frontend" sheet, inserted directly after Sheet1 -----------------
$frontend = $wb.Worksheets.Add([System.Type]::Missing, $sheet1)
$frontend.Name = "Frontend"

$frontend.Range("A2").Value = "React mobile web first"
$frontend.Range("A3").Value = "Pwa för att undvika appstore"
$frontend.Range("A4").Value = "Progressive web app"

# Leave the new sheet active, with A5 selected (the row right after the
# last bit of text) -- this also makes it the saved "activeTab".
$frontend.Range("A5").Select() | Out-Null
